$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9 updates
$ws.Cells.Item(9, 6).Value = 21492   # F9
$ws.Cells.Item(9, 7).Value = 19652   # G9
$ws.Cells.Item(9, 13).Value = 15353  # M9
$ws.Cells.Item(9, 14).Value = 11599  # N9
$ws.Cells.Item(9, 23).Value = 28087  # W9
$ws.Cells.Item(9, 24).Value = 13725  # X9

# Row 26 updates
$ws.Cells.Item(26, 3).Value = 0      # C26
$ws.Cells.Item(26, 4).Value = 0      # D26
$ws.Cells.Item(26, 5).Value = 0      # E26
$ws.Cells.Item(26, 6).Value = 3931   # F26
$ws.Cells.Item(26, 7).Value = 3663   # G26
$ws.Cells.Item(26, 13).Value = 2297  # M26
$ws.Cells.Item(26, 14).Value = 1782  # N26
$ws.Cells.Item(26, 23).Value = 2795  # W26
$ws.Cells.Item(26, 24).Value = 2041  # X26

# Row 31: mark H31 with the same "Bad" highlight formatting already used
# on H10/H11/H12/H28/H30 in column H. Copy the format from H30 so the
# underlying style record (border etc.) matches exactly.
$ws.Cells.Item(30, 8).Copy()
$ws.Cells.Item(31, 8).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Application.CalculateFullRebuild()
